$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Construction Site", "Block of Land" and "Fence" rows (original rows 9-11)
$ws.Rows("9:11").Delete() | Out-Null

# Remove the "Dump Truck" and "Spotlights" rows
# (originally rows 14-15, now rows 11-12 after the deletion above)
$ws.Rows("11:12").Delete() | Out-Null

# Rename "Wooden Plank" -> "Metal Beam" (now row 10)
$ws.Range("B10").Value = "Metal Beam"

# Mark the remaining construction items as "Done" and right align the status,
# matching the formatting used by the rest of the Status column
$ws.Range("G9").Value = "Done"
$ws.Range("G9").HorizontalAlignment = -4152

$ws.Range("G10").Value = "Done"
$ws.Range("G10").HorizontalAlignment = -4152

# Update the sheet's current selection to match the new data range
$ws.Range("H13:I19").Select() | Out-Null
